$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "massive MDY site template update": split the single "Date Sampled" date
# into separate Month / Day / Year columns, inserted right after "Transect"
# (old column D) and before the old "Date Sampled" column (old column E).
$ws.Range("E1:G1").EntireColumn.Insert()

# New column widths for E:G match the Site/Transect columns (C:D).
$ws.Range("E1:G1").EntireColumn.ColumnWidth = 8.67

# Carry the column header/data formatting from column D onto the new E:G
# columns so the style indices line up with the rest of the template.
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)

$rows = 2,3,4,5
foreach ($r in $rows) {
    $ws.Range("D" + $r).Copy()
    $ws.Range("E" + $r + ":G" + $r).PasteSpecial(-4122)
}

# Header row, written Day/Month/Year order so the shared-string table picks
# up the same insertion order as the source workbook.
$ws.Range("F1").Value = "Day"
$ws.Range("E1").Value = "Month"
$ws.Range("G1").Value = "Year"

# Populate Month/Day/Year for each data row from the existing "Date Sampled"
# value (now shifted right into column H), freezing the formula result down
# to a plain numeric literal.
foreach ($r in $rows) {
    $dateCell = "H" + $r

    $ws.Range("Z1").Formula = "=MONTH(" + $dateCell + ")"
    $monthVal = $ws.Range("Z1").Value2()

    $ws.Range("Z1").Formula = "=DAY(" + $dateCell + ")"
    $dayVal = $ws.Range("Z1").Value2()

    $ws.Range("Z1").Formula = "=YEAR(" + $dateCell + ")"
    $yearVal = $ws.Range("Z1").Value2()

    $ws.Range("Z1").Clear()

    $ws.Range("E" + $r).Value = $monthVal
    $ws.Range("F" + $r).Value = $dayVal
    $ws.Range("G" + $r).Value = $yearVal
}

# Restore the selection shown in the refreshed template.
$ws.Range("E1:G1").Select()
